$wb = $excel.ActiveWorkbook

# --- Sheet "Backtracking": mark rows 9 and 10 (Sudoku Solver, Word Search) Done on 29 Dec 2024 ---
$ws = $wb.Worksheets.Item("Backtracking")
$ws.Range("D9").Value = "Done"
$ws.Range("E9").Value = 45655
$ws.Range("D10").Value = "Done"
$ws.Range("E10").Value = 45655
$ws.Range("D12").Select()

# --- Sheet "Dynamic Programming": scroll/selection update ---
$ws2 = $wb.Worksheets.Item("Dynamic Programming")
$ws2.Activate()
$ws2.Range("D9").Select()

# --- Sheet "Divide and Conquer": selection update ---
$ws10 = $wb.Worksheets.Item("Divide and Conquer")
$ws10.Activate()
$ws10.Range("D4").Select()

$ws.Activate()
